$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.477.06'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '1.900.72'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').Value = "'245.02"
$ws.Range('E5').Value = '  +3.83%  '
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('D9').Value = "'0.338"
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').Value = "'0.0997"
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '2.180.61'
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('D13').Value = "'12.32"
$ws.Range('E13').Value = '  +7.93%  '
$ws.Range('D14').Value = '1.959.89'
$ws.Range('E14').Value = '  +5.84%  '
$ws.Range('D15').Value = "'0.688"
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = "'4.85"
$ws.Range('E16').Value = '  +3.19%  '
$ws.Range('D17').Value = '35.452.17'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = "'71.87"
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('D19').Value = '0.0₃0820'
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').Value = "'242.77"
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = "'171.91"
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = "'2.15"
$ws.Range('E26').Value = '  +15.20%  '
$ws.Range('D27').Value = "'8.54"
$ws.Range('E27').Value = '  +7.99%  '
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +23.39%  '
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').Value = "'4.15"
$ws.Range('E34').Value = '  +4.18%  '
$ws.Range('E35').Value = '  +7.38%  '
$ws.Range('E36').Value = '  +10.25%  '
$ws.Range('D37').Value = "'2.01"
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('E38').Value = '  +2.43%  '
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').Value = "'90.42"
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').Value = "'15.77"
$ws.Range('E41').Value = '  +4.62%  '
$ws.Range('D42').Value = '1.348.06'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').Value = "'0.0604"
$ws.Range('E43').Value = '  +12.37%  '
$ws.Range('D44').Value = "'48.72"
$ws.Range('E44').Value = '  +39.66%  '
$ws.Range('B45').Value = 'Gas'
$ws.Range('C45').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D45').Value = "'13.09"
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = "'2.34"
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Value = "'6.62"
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').Value = '2.087.11'
$ws.Range('E50').Value = '  +2.67%  '
$ws.Range('E51').Value = '  +1.59%  '
